# Add a TC for valid login
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1: header labels
$ws.Range("A1").Value = "login"
$ws.Range("B1").Value = "password"

# Row 2: existing valid-looking TC, keep values
$ws.Range("A2").Value = "test@test.com"
$ws.Range("B2").Value = "test123"

# Row 3: new TC data
$ws.Range("A3").Value = "adam1"
$ws.Range("B3").Value = "test123"

# Row 4: new TC data
$ws.Range("A4").Value = "asdfasdf"
$ws.Range("B4").Value = "catlover1"

# Row 5: new valid login test case - numeric login, left aligned
$ws.Range("A5").Value = 4445555
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("B5").Value = "ilovepasta5"

# Move the active selection
$ws.Range("D6").Select()
